# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet named "2022-Q3" right before the existing
#    "2022-Q2" sheet (so the final order is 总计, 2022-Q3, 2022-Q2, 2021-Q4,
#    2021-Q3), and populate it with the quarterly fund-holding detail table.
# 2. Insert a new row into the "总计" (summary) sheet so the 2022-Q3 totals
#    show up right after the header row, pushing the other quarters down.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")
$q2sheet = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------
# 1) Add the "2022-Q3" detail sheet, inserted before "2022-Q2".
# ---------------------------------------------------------------------
$q3sheet = $wb.Worksheets.Add($q2sheet)
$q3sheet.Name = "2022-Q3"

# Header row (bold / bordered "s=2" style is carried over implicitly by
# copying formatting from the existing 2022-Q2 sheet's header row).
$headerSrc = $q2sheet.Range("B1:H1")
$headerSrc.Copy()
$q3sheet.Range("B1").PasteSpecial()

$q3sheet.Range("B1").Value = "基金代码"
$q3sheet.Range("C1").Value = "基金名称"
$q3sheet.Range("D1").Value = "基金规模"
$q3sheet.Range("E1").Value = "股票总仓位"
$q3sheet.Range("F1").Value = "仓位占比"
$q3sheet.Range("G1").Value = "持有市值(亿元)"
$q3sheet.Range("H1").Value = "仓位排名"

# Index column (A) uses the same bold/bordered style as the existing sheets.
$aStyleSrc = $q2sheet.Range("A2")
$aStyleSrc.Copy()
$q3sheet.Range("A2:A8").PasteSpecial()

# Columns B, D, E, F, G hold numeric-looking figures that must stay TEXT
# (fund codes keep leading zeros, percentages/sizes keep trailing zeros) -
# format as text before writing so Excel doesn't coerce them to numbers.
$q3sheet.Range("B2:B8").NumberFormat = "@"
$q3sheet.Range("D2:G8").NumberFormat = "@"

$data = @(
    @("010861", "长信企业优选一年持有期灵活配置混合", "8.09", "81.28", "2.52", "0.2039", 10),
    @("014938", "同泰产业升级混合A",                  "1.01", "61.58", "2.92", "0.0295", 6),
    @("002409", "华夏新活力灵活配置混合A",              "0.15", "69.89", "5.68", "0.0085", 2),
    @("011361", "华夏博锐一年持有混合（MOM）A",          "0.11", "33.21", "4.71", "0.0052", 1),
    @("014939", "同泰产业升级混合C",                  "0.00", "61.58", "2.92", $null,     6),
    @("011362", "华夏博锐一年持有混合（MOM）C",          "0.00", "33.21", "4.71", $null,     1),
    @("002410", "华夏新活力灵活配置混合C",              "0.00", "69.89", "5.68", $null,     2)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]

    $q3sheet.Cells.Item($row, 1).Value = $i
    $q3sheet.Cells.Item($row, 2).Value = $rec[0]
    $q3sheet.Cells.Item($row, 3).Value = $rec[1]
    $q3sheet.Cells.Item($row, 4).Value = $rec[2]
    $q3sheet.Cells.Item($row, 5).Value = $rec[3]
    $q3sheet.Cells.Item($row, 6).Value = $rec[4]
    if ($rec[5] -eq $null) {
        $q3sheet.Cells.Item($row, 7).NumberFormat = "General"
        $q3sheet.Cells.Item($row, 7).Value = 0
    } else {
        $q3sheet.Cells.Item($row, 7).Value = $rec[5]
    }
    $q3sheet.Cells.Item($row, 8).Value = $rec[6]
}

# ---------------------------------------------------------------------
# 2) Insert the 2022-Q3 summary row into the "总计" sheet (row 2), pushing
#    the existing 2022-Q2 / 2021-Q4 / 2021-Q3 rows down by one.
# ---------------------------------------------------------------------
$summary.Rows.Item(2).Insert()

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 7
$summary.Cells.Item(2, 4).Value = 0.25

# Give the new A2 index cell the same bold/bordered style as the others.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
